# Generate Report for Handback
# Applies handback info (target file / handback file / handback datetime)
# to the zh-cn and de-de sheets, widens a few columns, and flips the
# Overview sheet's status text from "Ready for handoff" to
# "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$zhUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f55fb739a2becc0e49d35fbeb97cd5c857fdeee/e2e/90ae4d30-adb9-4c90-a532-eb17ba8b0997.md"

# ---- Overview sheet ----------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

$zhcn.Range("I2").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.md"
$zhcn.Range("J2").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.bdc58107db3ad851ca5abdf44805a75182a73397.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-03 11:07:26"

$zhcn.Range("I3").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.md"
$zhcn.Range("J3").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.bdc58107db3ad851ca5abdf44805a75182a73397.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-03 11:07:26"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $zhUrl, "", "", "90ae4d30-adb9-4c90-a532-eb17ba8b0997.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $zhUrl, "", "", "90ae4d30-adb9-4c90-a532-eb17ba8b0997.md") | Out-Null

# ---- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40

$dede.Range("I2").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.md"
$dede.Range("J2").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.bdc58107db3ad851ca5abdf44805a75182a73397.de-de.xlf"
$dede.Range("K2").Value = "2016-09-03 11:07:33"

$dede.Range("I3").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.md"
$dede.Range("J3").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.bdc58107db3ad851ca5abdf44805a75182a73397.de-de.xlf"
$dede.Range("K3").Value = "2016-09-03 11:07:33"

$dede.Hyperlinks.Add($dede.Range("I2"), $zhUrl, "", "", "90ae4d30-adb9-4c90-a532-eb17ba8b0997.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), $zhUrl, "", "", "90ae4d30-adb9-4c90-a532-eb17ba8b0997.md") | Out-Null
